$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# --- Weapons stats table (rows 1-5): add a new "Durability" row ---
# Insert a new row at row 6, pushing the rest of the sheet down by one.
$ws.Rows("6:6").Insert()

# The new row 6 becomes the new bottom of the table, so it should carry the
# thick-bottom-border formatting that used to belong to the old row 5 ("Cost").
$ws.Range("A5:F5").Copy()
$ws.Range("A6:F6").PasteSpecial(-4122)

# The old row 5 is now an interior row, so it should look like row 4 (no thick
# bottom border) instead.
$ws.Range("A4:F4").Copy()
$ws.Range("A5:F5").PasteSpecial(-4122)

# Label the new row.
$ws.Range("A6").Value = "Durability"

# --- Armor stats table: fill in the previously-blank "Defense" stat row ---
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 2
$ws.Range("D9").Value = 3

# Restore the selection to where it ends up after the edit.
[void]$ws.Range("C14").Select()
